$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Regenerate report: refresh handoff timestamps and update status for the
# 9ca4079f-3422-48eb-9e25-292383c9a0e4 file, which is now ready for handoff.

# Overview sheet
$overview.Range("D2").Value = "2016-03-18 17:28:22"
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-18 17:28:22"

# zh-cn sheet
$zhcn.Range("E2").Value = "2016-03-18 17:28:12"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-18 17:28:12"

# de-de sheet
$dede.Range("E2").Value = "2016-03-18 17:28:22"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-18 17:28:22"
